$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(6, 1).Value = "login-functionality;login-with-valid-username-and-password"
$ws.Cells.Item(6, 2).Value = "passed"
$ws.Cells.Item(6, 3).Value = "2020-12-12 12_05_13"
$ws.Cells.Item(6, 4).Value = "chrome"

$ws.Cells.Item(7, 1).Value = "nationalities-functionality;nationalities-all-functionality"
$ws.Cells.Item(7, 2).Value = "passed"
$ws.Cells.Item(7, 3).Value = "2020-12-12 12_05_48"
$ws.Cells.Item(7, 4).Value = "chrome"
